# "Generate Report for Handoff" - refresh the localization-status report:
#  - Status moves from "Handed back: in sync with en-US" to "Ready for handoff"
#  - the handoff timestamps advance to the new generation run
#  - the Status column narrows now that the new text is shorter than before

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$ws_overview.Range("E2").Value = "Ready for handoff"
$ws_overview.Range("F2").Value = "Ready for handoff"
$ws_overview.Range("G2").Value = "2016-08-29 11:00:45"

# --- zh-cn sheet ---
$ws_zhcn.Range("C2").Value = "Ready for handoff"
$ws_zhcn.Range("H2").Value = "2016-08-29 11:00:40"

# --- de-de sheet ---
$ws_dede.Range("C2").Value = "Ready for handoff"
$ws_dede.Range("H2").Value = "2016-08-29 11:00:45"

# --- Narrow the Status columns to fit the shorter text ---
# (engine quantizes ColumnWidth to 1/6-character steps; 16.3333... is the
# closest input that lands on the stored width nearest the new best-fit value)
$ws_overview.Range("E1:F1").ColumnWidth = 16.333333333333336
$ws_zhcn.Range("C1").ColumnWidth = 16.333333333333336
$ws_dede.Range("C1").ColumnWidth = 16.333333333333336
